# Apply convergence-table synthesis update: 23 data rows (was 18), refreshed
# shared strings and numeric results. Generated to match the target OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 23,9
$data[0,0] = 9
$data[0,1] = "Poisson"
$data[0,2] = "P1 FE"
$data[0,3] = 2
$data[0,4] = "Unstructured_triangles"
$data[0,5] = [double]"2.11982946328897"
$data[0,6] = "Triangles"
$data[0,7] = "Green"
$data[0,8] = [double]"31.36970090866089"
$data[1,0] = 17
$data[1,1] = "Poisson"
$data[1,2] = "P1 FE"
$data[1,3] = 3
$data[1,4] = "Unstructured_tetrahedra"
$data[1,5] = [double]"0.9881680535318813"
$data[1,6] = "Tetrahedra"
$data[1,7] = "Green"
$data[1,8] = [double]"166.0481879711151"
$data[2,0] = 21
$data[2,1] = "Poisson"
$data[2,2] = "VF5"
$data[2,3] = 2
$data[2,4] = "Deformed_quadrangles"
$data[2,5] = [double]"1.099715148758559"
$data[2,6] = "Squares"
$data[2,7] = "Green"
$data[2,8] = [double]"8.572041034698486"
$data[3,0] = 14
$data[3,1] = "Poisson"
$data[3,2] = "VF5"
$data[3,3] = 2
$data[3,4] = "Non_conforming_checkerboard"
$data[3,5] = [double]"0.3500538307526957"
$data[3,6] = "Squares"
$data[3,7] = "Green"
$data[3,8] = [double]"6.620848894119263"
$data[4,0] = 8
$data[4,1] = "Poisson"
$data[4,2] = "VF5"
$data[4,3] = 2
$data[4,4] = "Non_conforming_locally_refined"
$data[4,5] = [double]"0.9362211106984839"
$data[4,6] = "Squares"
$data[4,7] = "Green"
$data[4,8] = [double]"21.06601786613464"
$data[5,0] = 1
$data[5,1] = "Poisson"
$data[5,2] = "VF5"
$data[5,3] = 2
$data[5,4] = "Regular_brickwall"
$data[5,5] = [double]"-0.2073891905776415"
$data[5,6] = "Squares"
$data[5,7] = "Green"
$data[5,8] = [double]"2.427170038223267"
$data[6,0] = 3
$data[6,1] = "Poisson"
$data[6,2] = "VF5"
$data[6,3] = 2
$data[6,4] = "Regular_hexagons"
$data[6,5] = [double]"1.94163703307054"
$data[6,6] = "Hexagons"
$data[6,7] = "Green"
$data[6,8] = [double]"2.972708940505981"
$data[7,0] = 11
$data[7,1] = "Poisson"
$data[7,2] = "VF5"
$data[7,3] = 2
$data[7,4] = "Regular_squares"
$data[7,5] = [double]"2.009991317806771"
$data[7,6] = "Squares"
$data[7,7] = "Green"
$data[7,8] = [double]"3.357578992843628"
$data[8,0] = 15
$data[8,1] = "Poisson"
$data[8,2] = "VF5"
$data[8,3] = 2
$data[8,4] = "Unstructured_triangles"
$data[8,5] = [double]"0.6137798580984064"
$data[8,6] = "Triangles"
$data[8,7] = "Green"
$data[8,8] = [double]"3.177055835723877"
$data[9,0] = 20
$data[9,1] = "Poisson"
$data[9,2] = "VF5"
$data[9,3] = 3
$data[9,4] = "Non_conforming_checkerboard"
$data[9,5] = [double]"-0.2601624714817707"
$data[9,6] = "Cubes"
$data[9,7] = "Orange, BC violated. PB with mesh ?"
$data[9,8] = [double]"13.95365500450134"
$data[10,0] = 10
$data[10,1] = "Poisson"
$data[10,2] = "VF5"
$data[10,3] = 3
$data[10,4] = "Regular_cubes"
$data[10,5] = [double]"2.002870871514143"
$data[10,6] = "Cubes"
$data[10,7] = "Green"
$data[10,8] = [double]"24.51998209953308"
$data[11,0] = 19
$data[11,1] = "Poisson"
$data[11,2] = "VF5"
$data[11,3] = 3
$data[11,4] = "Unstructured_tetrahedra"
$data[11,5] = [double]"0.8110604778863585"
$data[11,6] = "Tetrahedra"
$data[11,7] = "Green"
$data[11,8] = [double]"37.72716212272644"
$data[12,0] = 5
$data[12,1] = "Poisson-Beltrami"
$data[12,2] = "P1 FE"
$data[12,3] = 2
$data[12,4] = "Unstructured_3D_triangles"
$data[12,5] = [double]"0.6091438624653882"
$data[12,6] = "3DTriangles"
$data[12,7] = "Green"
$data[12,8] = [double]"9.205935001373291"
$data[13,0] = 12
$data[13,1] = "Wave system"
$data[13,2] = "PStag scaling"
$data[13,3] = 2
$data[13,4] = "Regular squares"
$data[13,5] = [double]"-0.0002403283177935742"
$data[13,6] = "Squares"
$data[13,7] = "Green"
$data[13,8] = [double]"48.67472195625305"
$data[14,0] = 16
$data[14,1] = "Wave system"
$data[14,2] = "PStag scaling"
$data[14,3] = 2
$data[14,4] = "Unstructured triangles"
$data[14,5] = [double]"8.01376136056213e-05"
$data[14,6] = "Triangles"
$data[14,7] = "Orange"
$data[14,8] = [double]"138.2427699565887"
$data[15,0] = 0
$data[15,1] = "Wave system"
$data[15,2] = "Upwind"
$data[15,3] = 2
$data[15,4] = "Deformed quadrangles"
$data[15,5] = [double]"-3.487620212460162e-09"
$data[15,6] = "Deformed quadrangles"
$data[15,7] = "Green"
$data[15,8] = [double]"3.91550087928772"
$data[16,0] = 2
$data[16,1] = "Wave system"
$data[16,2] = "Upwind"
$data[16,3] = 2
$data[16,4] = "Regular brick wall"
$data[16,5] = [double]"-7.023350625076041e-06"
$data[16,6] = "Squares"
$data[16,7] = "Green"
$data[16,8] = [double]"6.973340034484863"
$data[17,0] = 22
$data[17,1] = "Wave system"
$data[17,2] = "Upwind"
$data[17,3] = 2
$data[17,4] = "Regular checkerboard"
$data[17,5] = [double]"-1.855096236805034e-11"
$data[17,6] = "Squares"
$data[17,7] = "Green"
$data[17,8] = [double]"6.204435110092163"
$data[18,0] = 18
$data[18,1] = "Wave system"
$data[18,2] = "Upwind"
$data[18,3] = 2
$data[18,4] = "Regular hexagons"
$data[18,5] = [double]"-8.564714050813453e-06"
$data[18,6] = "Hexagons"
$data[18,7] = "Green"
$data[18,8] = [double]"7.821051836013794"
$data[19,0] = 6
$data[19,1] = "Wave system"
$data[19,2] = "Upwind"
$data[19,3] = 2
$data[19,4] = "Regular squares"
$data[19,5] = [double]"-3.864405473494067e-05"
$data[19,6] = "Squares"
$data[19,7] = "Green"
$data[19,8] = [double]"6.330248832702637"
$data[20,0] = 4
$data[20,1] = "Wave system"
$data[20,2] = "Upwind"
$data[20,3] = 2
$data[20,4] = "Unstructured triangles"
$data[20,5] = [double]"-8.25752331143273e-11"
$data[20,6] = "Triangles"
$data[20,7] = "Green"
$data[20,8] = [double]"4.3265061378479"
$data[21,0] = 13
$data[21,1] = "Wave system"
$data[21,2] = "Upwind"
$data[21,3] = 3
$data[21,4] = "Regular cubes"
$data[21,5] = [double]"-3.26265785023925e-05"
$data[21,6] = "Cubes"
$data[21,7] = "Green"
$data[21,8] = [double]"11.58602809906006"
$data[22,0] = 7
$data[22,1] = "Wave system"
$data[22,2] = "Upwind"
$data[22,3] = 3
$data[22,4] = "Regular tetrahedra"
$data[22,5] = [double]"-1.117834915917594e-11"
$data[22,6] = "Tetrahedra"
$data[22,7] = "Green"
$data[22,8] = [double]"74.7223608493805"

$ws.Range("A2:I24").Value = $data

# Rows 19-24 are brand new in sheet1 (sheet previously ended at row 18); the
# A-column header style ("s=1": bold + thin border + center/top alignment)
# needs to be copied onto them explicitly so the whole A2:A24 run matches.
$ws.Range("A2").Copy()
$ws.Range("A19:A24").PasteSpecial(-4122)
$excel.CutCopyMode = 0
